$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix duplicated "类" typo in header row 1
$ws.Range("B1").Value = "其他食品类农村居民消费价格指数(上年=100)"
$ws.Range("D1").Value = "奶类农村居民消费价格指数(上年=100)"
$ws.Range("E1").Value = "干鲜瓜果类农村居民消费价格指数(上年=100)"
$ws.Range("H1").Value = "畜肉类农村居民消费价格指数(上年=100)"
$ws.Range("I1").Value = "禽肉类农村居民消费价格指数(上年=100)"
$ws.Range("K1").Value = "糖果糕点类农村居民消费价格指数(上年=100)"
$ws.Range("N1").Value = "薯类农村居民消费价格指数(上年=100)"
$ws.Range("O1").Value = "蛋类农村居民消费价格指数(上年=100)"
$ws.Range("Q1").Value = "豆类农村居民消费价格指数(上年=100)"

# Copy the year-label formatting (bold, centered, bordered) from A2 down to the new rows
$ws.Range("A2").Copy()
$ws.Range("A7:A8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Add new row 7: 2021年
$ws.Range("A7").Value = "2021年"
$ws.Range("B7").Value = 100.4
$ws.Range("C7").Value = 101.6
$ws.Range("D7").Value = 101.2
$ws.Range("E7").Value = 101.9
$ws.Range("F7").Value = 111
$ws.Range("G7").Value = 101.3
$ws.Range("H7").Value = 80.2
$ws.Range("I7").Value = 97.7
$ws.Range("J7").Value = 101.2
$ws.Range("K7").Value = 100.7
$ws.Range("L7").Value = 100.9
$ws.Range("M7").Value = 105
$ws.Range("N7").Value = 100.5
$ws.Range("O7").Value = 112.4
$ws.Range("P7").Value = 101.1
$ws.Range("Q7").Value = 106.5
$ws.Range("R7").Value = 98.8
$ws.Range("S7").Value = 97.7
$ws.Range("T7").Value = 105.9
$ws.Range("U7").Value = 102.7
$ws.Range("V7").Value = 105.6

# Add new row 8: 2022年 (only A8 and R8 populated)
$ws.Range("A8").Value = "2022年"
$ws.Range("R8").Value = 102.1
